$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# The application now grants access to a new table
# (integracao.tb_c_acesso_transac_integracao) instead of the old
# vw_bmh_online view for the block of rows 100-151.
$ws.Range("A100:A151").Value = "GRANT SELECT on integracao.tb_c_acesso_transac_integracao TO "

# Reflect the author's final on-screen state: scrolled/selected the
# newly edited D column cells.
$ws.Activate()
[void]$ws.Range("D100:D151").Select()
